$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 122, shifting existing rows 122-164 down to 124-166
$ws.Rows("122:123").Insert()

# Fill new row 122
$ws.Range("A122").Value = 1
$ws.Range("B122").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C122").Value = "Arica y Parinacota"
$ws.Range("D122").Value = 44900
$ws.Range("D122").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E122").Value = 15
$ws.Range("F122").Value = "Fruta"
$ws.Range("G122").Value = 100108
$ws.Range("H122").Value = "Tropicales y subtropicales"
$ws.Range("I122").Value = 100108002
$ws.Range("J122").Value = "Mango"
$ws.Range("K122").Value = "Sin especificar"
$ws.Range("L122").Value = "Especial"
$ws.Range("M122").Value = 780
$ws.Range("N122").Value = 5000
$ws.Range("O122").Value = 5500
$ws.Range("P122").Value = 5250
$ws.Range("Q122").Value = "`$/bandeja 4 kilos"
$ws.Range("R122").Value = "Perú"
$ws.Range("S122").Value = 1312
$ws.Range("T122").Value = 4

# Fill new row 123
$ws.Range("A123").Value = 1
$ws.Range("B123").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C123").Value = "Arica y Parinacota"
$ws.Range("D123").Value = 44900
$ws.Range("D123").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E123").Value = 15
$ws.Range("F123").Value = "Fruta"
$ws.Range("G123").Value = 100108
$ws.Range("H123").Value = "Tropicales y subtropicales"
$ws.Range("I123").Value = 100108002
$ws.Range("J123").Value = "Mango"
$ws.Range("K123").Value = "Sin especificar"
$ws.Range("L123").Value = "Primera"
$ws.Range("M123").Value = 600
$ws.Range("N123").Value = 5000
$ws.Range("O123").Value = 5500
$ws.Range("P123").Value = 5250
$ws.Range("Q123").Value = "`$/bandeja 4 kilos"
$ws.Range("R123").Value = "Perú"
$ws.Range("S123").Value = 1312
$ws.Range("T123").Value = 4
